$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header fix: "EXP ERIMENT COMPLETED" -> "EXPERIMENT COMPLETED" ---
$ws.Range("F1").Value = "EXPERIMENT COMPLETED"

# --- 2. Revise existing comment text for subject 1 (row 2), subject 2 (row 3), subject 3 (row 4) ---
$ws.Range("H2").Value = "TASK: it was not clear that the shuffles were starting from zero at every trial, so at the beginning he was trying to use as few as possible -> underline it more during the tutorial with feedback + safari version 16.4"
$ws.Range("H3").Value = "TUTORIAL (I part): bad quality images during tutorial; TASK: thought that using more shuffles resulted in more difficult following trials"
$ws.Range("H4").Value = 'TUTORIAL (I part): when starting tutorial, she was expecting to see the red version first (it is the first time the black one is shown first)  + BROWSER COMPATIBILITY: problems visualizing the text (FillText) with Safari (version 14.1.1 ) + it''s written "right or left" arrow, but instead it should be "left or right" + she did not feel tired at all, could have done other 2/3 blocks'

# --- 3. Complete subject 4 (row 5): fill in end time, duration, comment ---
$ws.Range("E5").Value = 0.65
$ws.Range("F5").Value = 0.67638888888888893
$ws.Range("G5").Value = 10
$ws.Range("F4").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("H4").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("H5").Value = "INSTRUCTIONS: 1. indicate which is the last page when it is possible to you can go back, before moving on with tutorial; 2. Avoid repeating same senteces over and over, since it is possible to go back; 3. Convey the idea that the task will become very difficult, with some difficult examples in the tutorial.  TASK: communicate score of last trial of block. "
$ws.Rows.Item(5).RowHeight = 49.2

# --- 4. New subject 5 (row 6) ---
$ws.Range("B5:H5").Copy()
$ws.Range("B6:H6").PasteSpecial(-4122)
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 23
$ws.Range("D6").Value = "F"
$ws.Range("E6").Value = 0.42152777777777778
$ws.Range("F6").Value = 0.43611111111111112
$ws.Range("G6").Value = 9
$ws.Range("H6").Value = "INTRO PAGE: change title (New Version of...) + you will receive INSTRUCTIONS…  + in case you HAVE QUESTIONS regarding + AND THAT of the other volunteers. BROWSER compatibility: (Firefox -> informed consent page not fitting, is cut. Also score is cut during the task). TUTORIAL (II part): low quality of images. TASK: she was expecting to see the score after pressing spacebar, like in the tutorial (she suggested to make it more consistent in the two cases. Maybe by inverting the score increase - solution images). "
$ws.Rows.Item(6).RowHeight = 90

# --- 5. New subject 6 (row 7) ---
$ws.Range("B6:H6").Copy()
$ws.Range("B7:H7").PasteSpecial(-4122)
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 24
$ws.Range("D7").Value = "M"
$ws.Range("E7").Value = 0.45694444444444443
$ws.Range("F7").Value = 0.46388888888888885
$ws.Range("G7").Value = 5
$ws.Range("H7").Value = "TASK: it was not clear what shuffles meant, he was thinking that one space bar press randomized the graph, and another one was taking the graph back to the previous visualization. Was very fast, but could not handle one trial more. Maybe he misunderstood / did not read the instructions carefully enough."
$ws.Rows.Item(7).RowHeight = 64.2

# --- 6. Placeholder subject rows 7-10 (rows 8-11), only SUBJECT NUMBER filled ---
$ws.Range("B8").Value = 7
$ws.Rows.Item(8).RowHeight = 37.2
$ws.Range("B9").Value = 8
$ws.Rows.Item(9).RowHeight = 38.4
$ws.Range("B10").Value = 9
$ws.Rows.Item(10).RowHeight = 35.4
$ws.Range("B11").Value = 10
$ws.Rows.Item(11).RowHeight = 24

# --- 7. Column widths (narrower B-G, much wider H to fit the long comments) ---
$ws.Range("B1").EntireColumn.ColumnWidth = 9.102120535714286
$ws.Range("C1").EntireColumn.ColumnWidth = 5.191964285714286
$ws.Range("D1").EntireColumn.ColumnWidth = 5.738839285714286
$ws.Range("E1").EntireColumn.ColumnWidth = 7.375558035714286
$ws.Range("F1").EntireColumn.ColumnWidth = 8.285714285714286
$ws.Range("G1").EntireColumn.ColumnWidth = 7.375558035714286
$ws.Range("H1").EntireColumn.ColumnWidth = 64.19196428571429

# --- 8. View state: scrolled down, new selection ---
$excel.ActiveWindow.ScrollRow = 5
$ws.Range("C8").Select()
